$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose updated price is a plain decimal number (e.g. "0.9987").
# Excel would otherwise auto-detect these as numeric values; the source data keeps
# them as plain text (as in the rest of the Price column), so force a Text format
# before writing the value, same as typing into a pre-formatted "Text" cell.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.082.72'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.860.95'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = '312.76'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.5082'
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("D8").Value = '0.3887'
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").Value = '0.08181'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '41.55'
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.169'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.847.76'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '20.12'
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.183'
$ws.Range("E15").Value = '  -1.72%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '0.9984'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001092'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '90.47'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.06661'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '17.54'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.955'
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '28.120.98'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '11.01'
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.222'
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = '2.076.49'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = '159.05'
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("D28").Value = '20.56'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Value = '2.403'
$ws.Range("E29").Value = '  -4.04%  '
$ws.Range("D30").Value = '125.23'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").Value = '0.1041'
$ws.Range("E31").Value = '  -2.17%  '
$ws.Range("D32").Value = '1.033'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").Value = '5.818'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '3.594'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '9.296'
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").Value = '0.02412'
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").Value = '0.06509'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '0.2179'
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").Value = '0.6410'
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").Value = '1.236'
$ws.Range("E40").Value = '  +2.25%  '
$ws.Range("D41").Value = '1.169'
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("D42").Value = '4.928'
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").Value = '11.08'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").Value = '0.6033'
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("D45").Value = '13.03'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.273'
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.659'
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("D48").Value = '1.982'
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").Value = '1.199'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = '121.07'
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").Value = '0.06869'
$ws.Range("E51").Value = '  +0.69%  '
